# Apply the edit described by the diff:
# - Add a new column "I" of data (praat false alarm / praat true detection /
#   Syll-o-matic true detection figures) to both tables in the sheet
#   (rows 8-15 for the TIMIT table, rows 32-39 for the LibriSpeech table).
# - The highlighted rows (12 and 36) get the teal "highlight" formatting that
#   matches the rest of that row (same as style used in column B/C/D/E there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- TIMIT dataset table ----
$ws.Range("I8").Value = "praat false alarm"
$ws.Range("I9").Value = 0.113554578168733

$ws.Range("I11").Value = "praat true detection"

$ws.Range("B12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 0.523990403838465

$ws.Range("I14").Value = "Syll-o-matic true detection"
$ws.Range("I15").Value = 0.49640143942423

# ---- LibriSpeech dataset table ----
$ws.Range("I32").Value = "praat false alarm"
$ws.Range("I33").Value = 0.182639237330657

$ws.Range("I35").Value = "praat true detection"

$ws.Range("B36").Copy()
$ws.Range("I36").PasteSpecial(-4122)
$ws.Range("I36").Value = 0.574410436527848

$ws.Range("I38").Value = "Syll-o-matic true detection"
$ws.Range("I39").Value = 0.476166583040642

# Move the active selection to I11, matching the workbook's last saved cursor position.
[void]$ws.Range("I11").Select()
